$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.330.74'
$ws.Range("E2").Value = '  -2.03%  '
$ws.Range("D3").Value = '3.164.87'
$ws.Range("E3").Value = '  -6.97%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '206.18'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -8.05%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '609.84'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -6.88%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.377'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -11.85%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.665'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.34%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '3.164.52'
$ws.Range("E10").Value = '  -6.87%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.533'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -16.84%  '
$ws.Range("E12").Value = '  +3.83%  '
$ws.Range("E13").Value = '  -17.83%  '
$ws.Range("D14").Value = '3.749.31'
$ws.Range("E14").Value = '  -6.51%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '5.24'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -7.41%  '
$ws.Range("D16").Value = '87.119.81'
$ws.Range("E16").Value = '  -2.09%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '31.92'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -14.70%  '
$ws.Range("D18").Value = '3.164.88'
$ws.Range("E18").Value = '  -6.12%  '
$ws.Range("E19").Value = '  -3.73%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '13.36'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -11.41%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '414.32'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -12.38%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '8.44'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -13.77%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '5.06'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -12.12%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '5.15'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -8.81%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '11.90'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -7.91%  '
$ws.Range("D26").Value = '3.336.04'
$ws.Range("E26").Value = '  -6.61%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '73.36'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -10.15%  '
$ws.Range("E28").Value = '  -12.54%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '0.158'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -18.22%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '541.66'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -8.96%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '8.18'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -14.64%  '
$ws.Range("E34").Value = '  -17.95%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '6.70'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -9.24%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '1.84'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -13.92%  '
$ws.Range("E37").Value = '  -10.13%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '21.73'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -10.28%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '21.79'
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '2.96'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -8.29%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '1.89'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -12.51%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.368'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -16.58%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '148.63'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.42%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '172.43'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -9.64%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '43.16'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -8.05%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '0.125'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.44%  '
$ws.Range("E49").Value = '  -15.70%  '
$ws.Range("E50").Value = '  -14.02%  '
$ws.Range("E51").Value = '  -13.06%  '
